$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Rules")

# Rule "R30" (row 10): Integer max (column C) was 18, corrected to 1.
$ws.Range("C10").Value = 1
